$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (008002057 / LUCIENE / 681500) - shifts all following rows up
$ws.Rows(2).Delete()
